# Generate Report for Archive
#
# The localization-status report moved every "Ready for handoff" cell to
# "In Translation" (Overview's per-locale roll-up columns, plus each
# locale sheet's own "Status" column). Since the new status text is
# shorter, the status columns are re-narrowed to fit it.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newStatusColumnWidth = 12.5   # narrower column width that fits "In Translation"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus) | Out-Null
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth   # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth   # de-de status column

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth       # Status column

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth       # Status column
